# Update "想去人数" (F column) figures that changed between data refreshes.
# Sheet "展览" (Exhibitions) row -> new F value
# Sheet "全部类型" (All types) row -> new F value

$wb = $excel.ActiveWorkbook

$updates1 = @{
    2  = 7092
    8  = 123
    12 = 205
    15 = 1834
    17 = 3670
    22 = 29
    23 = 2307
    32 = 1332
    33 = 118
}

$updates4 = @{
    2  = 7092
    9  = 123
    13 = 205
    16 = 1834
    18 = 3670
    23 = 29
    24 = 2307
    33 = 1332
    34 = 118
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
